# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, reflecting newly generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7775
    6  = 45
    9  = 6129
    10 = 165
    11 = 18
    12 = 34
    13 = 1848
    14 = 1381
    16 = 897
    17 = 185
    19 = 69
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
